$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all cell contents while preserving existing cell formatting
# (A1, B1, A11 currently carry the bold/border/centered header style).
$ws.Cells.ClearContents()

# --- Row 1: header ---
$ws.Range("A1").Value = "Gender"
$ws.Range("B1").Value = "Customer ID"

# --- Row 2 ---
$ws.Range("A2").Value = "Female"
$ws.Range("B2").Value = 2398

# --- Row 3 ---
$ws.Range("A3").Value = "Male"
$ws.Range("B3").Value = 5105

# --- Row 6: "Phân tích / Nhận xét" label, reusing the header style ---
$ws.Range("A1").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").Value = "Phân tích / Nhận xét"

# --- Row 7: analysis text ---
$text = "Here's a brief report based on the pie chart:`nThe pie chart displays the gender distribution of customers.  Males constitute a larger portion of the customer base at 68%, while females represent 32%. This suggests the product/service might be more appealing or marketed more effectively towards men.  Further investigation into customer demographics and marketing strategies could reveal opportunities for better engaging the female demographic."
$ws.Range("A7").Value = $text

# The old row 11 (former "Phân tích" cell) still carries the header style
# even though its contents were cleared; drop the formatting entirely so
# it doesn't inflate the sheet's used range.
$ws.Range("A11").Clear()

# Rename the sheet to match the new content.
$ws.Name = "Sheet 0"
